$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New issue row (#27): Billing screen issue reported by Sabethan ---
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "Billing"
$ws.Cells.Item(28, 3).Value = "Add first transaction without selecting the product (click '+' button)"
$ws.Cells.Item(28, 4).Value = "Click '+' button showing the confirm dialog box ""Do you want to delete the item from the list?"" and clicks cancel throws exception"
$ws.Cells.Item(28, 6).Value = "Sabethan"

# --- Existing issue row (#26, "Test class missing...") gets marked fixed ---
$ws.Cells.Item(27, 8).Value = "Sabethan"

# --- New "Fixed Remark" column header (matches the look of the other header cells) ---
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "Fixed Remark"
$ws.Cells.Item(1, 9).VerticalAlignment = -4160
$ws.Cells.Item(1, 9).WrapText = $false

$ws.Cells.Item(27, 9).Value = "Test class added for storekeeper"

# --- Move view/selection down to the newly edited rows ---
$ws.Rows.Item(19).Select() | Out-Null
